$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current header row (row 2), pushing the
# header row to row 4 and all the data rows down by two (rows 3-12 -> 5-14).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# B2 reuses the existing bold header style (copy format from the header
# row, which now lives at row 4).
$ws.Range("B4").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2: new "Autores:" / author name cells.
$ws.Range("B2").Value = "Autores:"
$ws.Range("C2").Value = "Hugo Grochau"
$ws.Range("D2").Value = "Leonardo Kaplan"

# C2 / D2 get a new style: thin box border all around + centered text.
$ws.Range("C2:D2").Borders.Weight = 2
$ws.Range("C2:D2").HorizontalAlignment = -4108

$ws.Range("E18").Select()
